# ============================================================================
# Edit script: "Add files via upload" commit
#  1. Slide 1 (sldId 260): two text tweaks on existing shapes (id 5 and id 9),
#     plus a resize/move of shape id 9.
#  2. A brand-new "Welcome Page" slide is inserted at position 2 (sldId 267).
#  3. Slide with sldId 263 (5th slide): a couple of label edits, a line/table
#     reflow, and a new "Collapse button" caption text box.
# ============================================================================

$emuPerPt = 12700.0
function Emu([double]$v) { return $v / $emuPerPt }

function Find-ShapeById($shapes, $targetId) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $targetId) {
            return $sh
        }
        if ($sh.Type -eq 6) {
            $found = Find-ShapeById $sh.GroupItems $targetId
            if ($found) { return $found }
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# 1. Slide 1 (sldId 260) -- rename "Front PAGE" -> "Home PAGE", and the
#    "Home" pill becomes "Dashboard" with a new position/size.
# ----------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

$sh5 = Find-ShapeById $slide1.Shapes 5
$sh5.TextFrame.TextRange.Text = "Home PAGE(App Name instead)"

$sh9 = Find-ShapeById $slide1.Shapes 9
$sh9.Left = Emu 4296834
$sh9.Top = Emu 755897
$sh9.Width = Emu 1143169
$sh9.Height = Emu 796509
$sh9.TextFrame.TextRange.Text = "Dashboard"

# ----------------------------------------------------------------------
# 3. Slide with sldId 263 -- grab it *before* the insert below changes
#    the positional index (it is the 5th slide prior to insertion).
# ----------------------------------------------------------------------
$slide263 = $p.Slides.Item(5)

$shHeader = Find-ShapeById $slide263.Shapes 5
$shHeader.TextFrame.TextRange.Text = "Dashboard(App icon)"

$shBg = Find-ShapeById $slide263.Shapes 4
$shBg.TextFrame.TextRange.Text = "C"

$shLine = Find-ShapeById $slide263.Shapes 21
$shLine.Left = Emu 1173099
$shLine.Top = Emu 1192289
$shLine.Width = Emu 10029444
$shLine.Height = Emu 0

$shTable = Find-ShapeById $slide263.Shapes 22
$shTable.Left = Emu 1741784
$shTable.Top = Emu 1593751
$shTable.Width = Emu 2388542
$shTable.Height = Emu 4232819

$newCaption = $slide263.Shapes.AddTextbox(1, (Emu 2226792), (Emu 1179051), (Emu 1700784), (Emu 307777))
$newCaption.TextFrame.WordWrap = -1
$newCaption.TextFrame.TextRange.Text = "Collapse button"
$newCaption.TextFrame.TextRange.Font.Size = 14

# ----------------------------------------------------------------------
# 2. Insert the brand-new "Welcome Page" slide at position 2 (blank layout).
# ----------------------------------------------------------------------
$newSlide = $p.Slides.Add(2, 12)

$tb = $newSlide.Shapes.AddTextbox(1, (Emu 5202936), (Emu 356616), (Emu 1566839), (Emu 369332))
$tb.TextFrame.TextRange.Text = "Welcome Page"

# -- Outer group (Group 2): the "card" that holds login controls --------
$rectBorder = $newSlide.Shapes.AddShape(1, (Emu 1380744), (Emu 292608), (Emu 10029443), (Emu 6382298))
$rectBorder.Name = "Rectangle 9"
$rectBorder.TextFrame.TextRange.Text = ""

$tbParen = $newSlide.Shapes.AddTextbox(1, (Emu 5715670), (Emu 5047114), (Emu 3719322), (Emu 420538))
$tbParen.Name = "TextBox 10"
$tbParen.TextFrame.TextRange.Text = ")"
$tbParen.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$tbImgBg = $newSlide.Shapes.AddTextbox(1, (Emu 8531352), (Emu 5782672), (Emu 2441448), (Emu 369332))
$tbImgBg.Name = "TextBox 11"
$tbImgBg.TextFrame.TextRange.Text = "Image in bg"

$welcomeBanner = $newSlide.Shapes.AddShape(5, (Emu 1945424), (Emu 641355), (Emu 9262872), (Emu 878778))
$welcomeBanner.Name = "Rectangle: Rounded Corners 13"
$welcomeBanner.TextFrame.TextRange.Text = "Welcome!  "

# Group the four shapes above into "Group 3"
$innerRange = $newSlide.Shapes.Range(@($rectBorder.Name, $tbParen.Name, $tbImgBg.Name, $welcomeBanner.Name))
$innerGroup = $innerRange.Group()
$innerGroup.Name = "Group 3"

$contactUs = $newSlide.Shapes.AddShape(5, (Emu 3194732), (Emu 5418442), (Emu 1249305), (Emu 813110))
$contactUs.Name = "Rectangle: Rounded Corners 4"
$contactUs.TextFrame.TextRange.Text = "Contact Us"

$news = $newSlide.Shapes.AddShape(5, (Emu 5441718), (Emu 5492610), (Emu 1249305), (Emu 813110))
$news.Name = "Rectangle: Rounded Corners 5"
$news.TextFrame.TextRange.Text = "News"

$logo = $newSlide.Shapes.AddShape(5, (Emu 5968805), (Emu 1729875), (Emu 1444436), (Emu 435753))
$logo.Name = "Rectangle: Rounded Corners 6"
$logo.TextFrame.TextRange.Text = "LOGO"

$connector = $newSlide.Shapes.AddConnector(1, (Emu 1380744), (Emu 3483757), (Emu 11410188), (Emu 3483757))
$connector.Name = "Straight Connector 7"
$connector.ConnectorFormat.BeginConnect($rectBorder, 1)
$connector.ConnectorFormat.EndConnect($rectBorder, 3)

$login = $newSlide.Shapes.AddShape(5, (Emu 2153921), (Emu 2634729), (Emu 2878836), (Emu 472349))
$login.Name = "Rectangle: Rounded Corners 8"
$login.TextFrame.TextRange.Text = "Login"

# Group "Group 3" plus the four standalone shapes/connector into "Group 2"
$outerRange = $newSlide.Shapes.Range(@($innerGroup.Name, $contactUs.Name, $news.Name, $logo.Name, $connector.Name, $login.Name))
$outerGroup = $outerRange.Group()
$outerGroup.Name = "Group 2"
$outerGroup.Left = Emu 1380744
$outerGroup.Top = Emu 822959
$outerGroup.Width = Emu 9236202
$outerGroup.Height = Emu 5605165

$register = $newSlide.Shapes.AddShape(5, (Emu 7650722), (Emu 2897037), (Emu 2651145), (Emu 414834))
$register.Name = "Rectangle: Rounded Corners 20"
$register.TextFrame.TextRange.Text = "Register"
